$wb = $excel.ActiveWorkbook

# --- Layer0 sheet ---
$ws0 = $wb.Worksheets.Item("Layer0")

$ws0.Range("B2").Value = -0.7048016320072106
$ws0.Range("C2").Value = 0.0006091205190841796

$ws0.Range("B3").Value = 0.4015493903758083
$ws0.Range("C3").Value = -0.5052744488970066

$ws0.Range("B4").Value = -1.384040654419918
$ws0.Range("C4").Value = -0.705614273029844

# --- Layer1 sheet ---
$ws1 = $wb.Worksheets.Item("Layer1")

$ws1.Range("B2").Value = -0.5909289811262698
$ws1.Range("C2").Value = -0.3652495545693971

$ws1.Range("B3").Value = -0.7273794528638384
$ws1.Range("C3").Value = 0.6263236006726336

$ws1.Range("B4").Value = -0.9386349892913847
$ws1.Range("C4").Value = 0.5491558014066598
